$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "valores" column header (shared string) in C1
$ws.Range("C1").Value = "valores"

# The 99 data values for C2:C100 (ages expressed as decimal years,
# computed elsewhere and pasted in as literal values)
$valores = @(
    34.5,
    35.5,
    44.5,
    48.1666666666667,
    53.1666666666667,
    58.1666666666667,
    63.1666666666667,
    68.1666666666667,
    73.1666666666667,
    78.1666666666668,
    83.1666666666668,
    88.1666666666668,
    93.1666666666668,
    98.1666666666668,
    103.166666666667,
    108.166666666667,
    113.166666666667,
    118.166666666667,
    123.166666666667,
    128.166666666667,
    133.166666666667,
    138.166666666667,
    143.166666666667,
    148.166666666667,
    153.166666666667,
    158.166666666667,
    163.166666666667,
    168.166666666667,
    173.166666666667,
    178.166666666667,
    183.166666666667,
    188.166666666667,
    193.166666666667,
    198.166666666667,
    203.166666666667,
    208.166666666667,
    213.166666666667,
    218.166666666667,
    223.166666666667,
    228.166666666667,
    233.166666666667,
    238.166666666667,
    243.166666666667,
    248.166666666667,
    253.166666666667,
    258.166666666667,
    263.166666666667,
    268.166666666667,
    273.166666666667,
    278.166666666667,
    283.166666666667,
    288.166666666667,
    293.166666666667,
    298.166666666667,
    303.166666666667,
    308.166666666667,
    313.166666666667,
    318.166666666667,
    323.166666666667,
    328.166666666667,
    333.166666666667,
    338.166666666667,
    343.166666666667,
    348.166666666667,
    353.166666666667,
    358.166666666667,
    363.166666666667,
    368.166666666667,
    373.166666666667,
    378.166666666667,
    383.166666666667,
    388.166666666667,
    393.166666666667,
    398.166666666667,
    403.166666666667,
    408.166666666667,
    413.166666666667,
    418.166666666667,
    423.166666666667,
    428.166666666667,
    433.166666666667,
    438.166666666667,
    443.166666666667,
    448.166666666667,
    453.166666666667,
    458.166666666667,
    463.166666666667,
    468.166666666667,
    473.166666666667,
    478.166666666667,
    483.166666666667,
    488.166666666667,
    493.166666666667,
    498.166666666667,
    503.166666666667,
    508.166666666667,
    513.166666666667,
    518.166666666667,
    523.166666666667
)

for ($i = 0; $i -lt $valores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $valores[$i]
}

# First three rows (C2:C4) keep the default/general number format;
# C5:C101 (including the trailing blank row) get a one-decimal format
$ws.Range("C5:C101").NumberFormat = "0.0"

# C101 is a formatted-but-empty trailing cell
$ws.Range("C101").ClearContents()

# Column C sizing to fit its new content (best-fit width ~12.57 chars,
# matching what Excel computes for this font/content)
$ws.Columns("C").ColumnWidth = 11.71

# Match the saved selection/active cell from the source workbook
[void]$ws.Range("C1").Select()
